# Apply updated crypto price/volume data as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "44.172.39", "1.01", "0.0(subscript3)0962")
# formatted with local-style dot separators. Excel's COM layer will
# silently coerce plain numeric-looking strings (e.g. "1.01") into real
# numbers when assigned via .Value, which would change the stored cell
# type from text to number. Force the number format to Text ("@") on
# each D cell we are about to rewrite so the value is kept as a string,
# matching the original workbook's inlineStr cells.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "44.172.39"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "2.238.89"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "306.38"
$ws.Range("E5").Value = "  -2.85%  "

$ws.Range("D6").Value = "93.95"
$ws.Range("E6").Value = "  -5.16%  "

$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("D10").Value = "34.51"
$ws.Range("E10").Value = "  -4.34%  "

$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").Value = "7.15"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.331.71"
$ws.Range("E14").Value = "  +3.98%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.580.50"
$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("D16").Value = "0.827"

$ws.Range("D17").Value = "13.48"
$ws.Range("E17").Value = "  -3.33%  "

$ws.Range("D18").Value = "43.927.90"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").Value = "0.0₃0962"
$ws.Range("E19").Value = "  -1.40%  "

$ws.Range("D20").Value = "6.36"
$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("D21").Value = "12.06"
$ws.Range("E21").Value = "  -8.21%  "

$ws.Range("D22").Value = "65.50"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "3.11"
$ws.Range("E23").Value = "  +4.35%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "237.50"
$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("D25").Value = "1.99"
$ws.Range("E25").Value = "  -1.35%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  -2.68%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  +3.56%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "38.07"
$ws.Range("E29").Value = "  +4.79%  "

$ws.Range("D30").Value = "19.99"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").Value = "5.82"
$ws.Range("E31").Value = "  -2.45%  "

$ws.Range("D32").Value = "152.86"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("D33").Value = "0.0792"
$ws.Range("E33").Value = "  -5.05%  "

$ws.Range("D34").Value = "2.61"
$ws.Range("E34").Value = "  -1.57%  "

$ws.Range("E35").Value = "  -4.90%  "

$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  +1.59%  "

$ws.Range("D37").Value = "0.106"
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -7.55%  "

$ws.Range("D39").Value = "3.49"
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("D40").Value = "3.80"
$ws.Range("E40").Value = "  -4.85%  "

$ws.Range("D41").Value = "14.33"
$ws.Range("E41").Value = "  -8.26%  "

$ws.Range("D42").Value = "0.0297"
$ws.Range("E42").Value = "  -3.07%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "1.741.95"
$ws.Range("E44").Value = "  +2.28%  "

$ws.Range("D45").Value = "82.35"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "0.190"
$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").Value = "14.95"
$ws.Range("E47").Value = "  +8.35%  "

$ws.Range("D48").Value = "99.44"
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("D49").Value = "4.91"
$ws.Range("E49").Value = "  -4.96%  "

$ws.Range("D50").Value = "8.04"
$ws.Range("E50").Value = "  -1.24%  "

$ws.Range("D51").Value = "1.56"
$ws.Range("E51").Value = "  -2.41%  "
